# Edit script: update Column I (BI_T2) values with corrected Barthel Index percentages,
# rename header D1 from "MRS_T0" to "mrs_T0", and clear 3 stray I-column values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header label casing for column D (MRS_T0 -> mrs_T0)
$ws.Range("D1").Value = "mrs_T0"

# Update column I (BI_T2) values with corrected figures
$ws.Cells.Item(4, 9).Value = 95
$ws.Cells.Item(7, 9).Value = 85
$ws.Cells.Item(10, 9).Value = 85
$ws.Cells.Item(16, 9).Value = 95
$ws.Cells.Item(19, 9).Value = 100
$ws.Cells.Item(22, 9).Value = 85
$ws.Cells.Item(25, 9).Value = 90
$ws.Cells.Item(30, 9).Value = 85
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(38, 9).Value = 15
$ws.Cells.Item(44, 9).Value = 40
$ws.Cells.Item(47, 9).Value = 95
$ws.Cells.Item(49, 9).Value = 25
$ws.Cells.Item(51, 9).Value = 95
$ws.Cells.Item(54, 9).Value = 25
$ws.Cells.Item(57, 9).Value = 100
$ws.Cells.Item(66, 9).Value = 65
$ws.Cells.Item(74, 9).Value = 85
$ws.Cells.Item(82, 9).Value = 100
$ws.Cells.Item(88, 9).Value = 100
$ws.Cells.Item(91, 9).Value = 100
$ws.Cells.Item(94, 9).Value = 20
$ws.Cells.Item(97, 9).Value = 100
$ws.Cells.Item(108, 9).Value = 100
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(114, 9).Value = 95
$ws.Cells.Item(123, 9).Value = 100
$ws.Cells.Item(129, 9).Value = 100
$ws.Cells.Item(135, 9).Value = 85
$ws.Cells.Item(141, 9).Value = 80
$ws.Cells.Item(143, 9).Value = 35
$ws.Cells.Item(146, 9).Value = 45
$ws.Cells.Item(152, 9).Value = 80
$ws.Cells.Item(161, 9).Value = 100
$ws.Cells.Item(164, 9).Value = 100
$ws.Cells.Item(167, 9).Value = 0
$ws.Cells.Item(169, 9).Value = 70
$ws.Cells.Item(172, 9).Value = 85
$ws.Cells.Item(184, 9).Value = 95
$ws.Cells.Item(186, 9).Value = 45
$ws.Cells.Item(189, 9).Value = 95
$ws.Cells.Item(192, 9).Value = 75
$ws.Cells.Item(197, 9).Value = 100
$ws.Cells.Item(200, 9).Value = 45
$ws.Cells.Item(203, 9).Value = 35
$ws.Cells.Item(216, 9).Value = 100
$ws.Cells.Item(219, 9).Value = 50
$ws.Cells.Item(222, 9).Value = 100
$ws.Cells.Item(225, 9).Value = 100
$ws.Cells.Item(228, 9).Value = 100
$ws.Cells.Item(237, 9).Value = 70
$ws.Cells.Item(240, 9).Value = 100
$ws.Cells.Item(243, 9).Value = 25
$ws.Cells.Item(246, 9).Value = 45
$ws.Cells.Item(249, 9).Value = 30
$ws.Cells.Item(257, 9).Value = 100
$ws.Cells.Item(262, 9).Value = 25

# Clear stray column I values that should not be present
$ws.Cells.Item(265, 9).ClearContents()
$ws.Cells.Item(277, 9).ClearContents()
$ws.Cells.Item(286, 9).ClearContents()
